$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2000
$ws.Range("J125").Value = 2000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920

$ws.Range("H129").Value = 10000000
$ws.Range("J129").Value = 10000000
$ws.Range("L129").Value = 30000000
$ws.Range("N129").Value = -30010000

$ws.Range("H132").Value = 4373.0527
$ws.Range("I132").Value = 4443
$ws.Range("K132").Value = 13329
$ws.Range("M132").Value = -10799

$ws.Range("H137").Value = 44950.13
$ws.Range("I137").Value = 1253.8572
$ws.Range("K137").Value = 3761.5716
$ws.Range("M137").Value = -1211.5716

$ws.Range("H138").Value = 1833.5
$ws.Range("I138").Value = 1233.2333
$ws.Range("J138").Value = 2307.3948
$ws.Range("K138").Value = 3699.699900000001
$ws.Range("L138").Value = 6922.1844
$ws.Range("M138").Value = 1440.300099999999
$ws.Range("N138").Value = -17202.1844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1641.9
$ws.Range("I2").Value = 1441.0714
$ws.Range("K2").Value = 1441.0714
$ws.Range("M2").Value = -1328.0714

$ws.Range("H16").Value = 640
$ws.Range("I16").Value = 640
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 640
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -353
$ws.Range("N16").ClearContents()

$ws.Range("H45").Value = 4010.8572
$ws.Range("I45").Value = 4049.3333
$ws.Range("J45").Value = 3982
$ws.Range("K45").Value = 4049.3333
$ws.Range("L45").Value = 3982
$ws.Range("M45").Value = -3672.3333
$ws.Range("N45").Value = -4736

$ws.Range("H74").Value = 41669436
$ws.Range("J74").Value = 2502.3333
$ws.Range("L74").Value = 2502.3333
$ws.Range("N74").Value = -4250.3333

$ws.Range("H77").Value = 41669436
$ws.Range("J77").Value = 2502.3333
$ws.Range("L77").Value = 12511.6665
$ws.Range("N77").Value = -21247.6665

$ws.Range("H116").Value = 1641.9
$ws.Range("I116").Value = 1441.0714
$ws.Range("K116").Value = 1441.0714
$ws.Range("M116").Value = 852.9286

$ws.Range("H122").Value = 1972.2354
$ws.Range("I122").Value = 1980.9333
$ws.Range("J122").Value = 1907
$ws.Range("K122").Value = 5942.7999
$ws.Range("L122").Value = 5721
$ws.Range("M122").Value = -3492.7999
$ws.Range("N122").Value = -10621

$ws.Range("H132").Value = 13220.884
$ws.Range("I132").Value = 1355.7333
$ws.Range("K132").Value = 4067.199900000001
$ws.Range("M132").Value = -1537.199900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1641.9
$ws.Range("I3").Value = 1441.0714
$ws.Range("K3").Value = 1441.0714
$ws.Range("M3").Value = -1327.0714

$ws.Range("H20").Value = 3022.5293
$ws.Range("I20").Value = 3981.7778
$ws.Range("J20").Value = 1943.375
$ws.Range("K20").Value = 3981.7778
$ws.Range("L20").Value = 1943.375
$ws.Range("M20").Value = -3734.7778
$ws.Range("N20").Value = -2437.375

$ws.Range("H94").Value = 864.94116
$ws.Range("I94").Value = 816.2381
$ws.Range("K94").Value = 816.2381
$ws.Range("M94").Value = -365.2381

$ws.Range("H99").Value = 1971.5385
$ws.Range("I99").Value = 1436.6666
$ws.Range("J99").Value = 2430
$ws.Range("K99").Value = 1436.6666
$ws.Range("L99").Value = 2430
$ws.Range("M99").Value = 61.33339999999998
$ws.Range("N99").Value = -5426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1181.1111
$ws.Range("J16").Value = 1125
$ws.Range("L16").Value = 1125
$ws.Range("N16").Value = -1699

$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 250

$ws.Range("H105").Value = 5953125.5
$ws.Range("I105").Value = 9615895
$ws.Range("K105").Value = 9615895
$ws.Range("M105").Value = -9614148

$ws.Range("H113").Value = 1181.1111
$ws.Range("J113").Value = 1125
$ws.Range("L113").Value = 1125
$ws.Range("N113").Value = -5465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 3880.7742
$ws.Range("J121").Value = 4542.269
$ws.Range("L121").Value = 13626.807
$ws.Range("N121").Value = -16246.807

$ws.Range("H131").Value = 713.52
$ws.Range("J131").Value = 714.6667
$ws.Range("L131").Value = 2144.0001
$ws.Range("N131").Value = -12224.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3686782.5
$ws.Range("I70").Value = 18957.143
$ws.Range("J70").Value = 6254260.5
$ws.Range("K70").Value = 18957.143
$ws.Range("L70").Value = 6254260.5
$ws.Range("M70").Value = -18687.143
$ws.Range("N70").Value = -6254800.5

$ws.Range("H73").Value = 3686782.5
$ws.Range("I73").Value = 18957.143
$ws.Range("J73").Value = 6254260.5
$ws.Range("K73").Value = 18957.143
$ws.Range("L73").Value = 6254260.5
$ws.Range("M73").Value = -18021.143
$ws.Range("N73").Value = -6256132.5

$ws.Range("H80").Value = 3447.0588
$ws.Range("J80").Value = 3560
$ws.Range("L80").Value = 3560
$ws.Range("N80").Value = -5556

$ws.Range("H83").Value = 3447.0588
$ws.Range("J83").Value = 3560
$ws.Range("L83").Value = 17800
$ws.Range("N83").Value = -27784

$ws.Range("H97").Value = 465
$ws.Range("I97").Value = 465
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 465
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 31
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 1134.1177
$ws.Range("I102").Value = 1190.0834
$ws.Range("J102").Value = 999.8
$ws.Range("K102").Value = 1190.0834
$ws.Range("L102").Value = 999.8
$ws.Range("M102").Value = 431.9166
$ws.Range("N102").Value = -4243.8

$ws.Range("H122").Value = 3150
$ws.Range("I122").Value = 2300
$ws.Range("K122").Value = 6900
$ws.Range("M122").Value = -4450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3258.12
$ws.Range("I7").Value = 3197.5264
$ws.Range("J7").Value = 3450
$ws.Range("K7").Value = 3197.5264
$ws.Range("L7").Value = 3450
$ws.Range("M7").Value = -3085.5264
$ws.Range("N7").Value = -3674

$ws.Range("H40").Value = 3880.652
$ws.Range("I40").Value = 3230.8
$ws.Range("J40").Value = 4061.1667
$ws.Range("K40").Value = 3230.8
$ws.Range("L40").Value = 4061.1667
$ws.Range("M40").Value = -3094.8
$ws.Range("N40").Value = -4333.1667

$ws.Range("H61").Value = 3463.08
$ws.Range("I61").Value = 1293.3158
$ws.Range("J61").Value = 10334
$ws.Range("K61").Value = 1293.3158
$ws.Range("L61").Value = 10334
$ws.Range("M61").Value = -1091.3158
$ws.Range("N61").Value = -10738

$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

$ws.Range("H113").Value = 3463.08
$ws.Range("I113").Value = 1293.3158
$ws.Range("J113").Value = 10334
$ws.Range("K113").Value = 1293.3158
$ws.Range("L113").Value = 10334
$ws.Range("M113").Value = 876.6841999999999
$ws.Range("N113").Value = -14674

$ws.Range("H126").Value = 3258.12
$ws.Range("I126").Value = 3197.5264
$ws.Range("J126").Value = 3450
$ws.Range("K126").Value = 9592.5792
$ws.Range("L126").Value = 10350
$ws.Range("M126").Value = -7122.5792
$ws.Range("N126").Value = -15290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 913
$ws.Range("I126").Value = 781.53845
$ws.Range("J126").Value = 1157.1428
$ws.Range("K126").Value = 2344.61535
$ws.Range("L126").Value = 3471.4284
$ws.Range("M126").Value = 125.38465
$ws.Range("N126").Value = -8411.428400000001

$ws.Range("H132").Value = 1262.2667
$ws.Range("I132").Value = 711.3333
$ws.Range("J132").Value = 3466
$ws.Range("K132").Value = 2133.9999
$ws.Range("L132").Value = 10398
$ws.Range("M132").Value = 396.0001000000002
$ws.Range("N132").Value = -15458
